$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Icam1"
$ws.Cells.Item(2,3).Value = "Itgb2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 163.9108346666667
$ws.Cells.Item(2,8).Value = 491.732504
$ws.Cells.Item(2,9).Value = 0.8426759240348239
$ws.Cells.Item(2,10).Value = 0.8426759240348242
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 83.95844533333333
$ws.Cells.Item(2,14).Value = 251.875336
$ws.Cells.Item(2,15).Value = 0.997975448786732
$ws.Cells.Item(2,16).Value = 0.997975448786732
$ws.Cells.Item(2,17).Value = 13761.69885190237
$ws.Cells.Item(2,18).Value = 123855.2896671213
$ws.Cells.Item(2,19).Value = 0.8409698834704274
$ws.Cells.Item(2,20).Value = 0.8409698834704277

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Icam1"
$ws.Cells.Item(3,3).Value = "Itgb2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 163.9108346666667
$ws.Cells.Item(3,8).Value = 491.732504
$ws.Cells.Item(3,9).Value = 0.8426759240348239
$ws.Cells.Item(3,10).Value = 0.8426759240348242
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.170323
$ws.Cells.Item(3,14).Value = 0.510969
$ws.Cells.Item(3,15).Value = 0.002024551213268089
$ws.Cells.Item(3,16).Value = 0.00202455121326809
$ws.Cells.Item(3,17).Value = 27.91778509293066
$ws.Cells.Item(3,18).Value = 251.260065836376
$ws.Cells.Item(3,19).Value = 0.001706040564396511
$ws.Cells.Item(3,20).Value = 0.001706040564396512

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Icam1"
$ws.Cells.Item(4,3).Value = "Itgb2"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 29.35342966666667
$ws.Cells.Item(4,8).Value = 88.060289
$ws.Cells.Item(4,9).Value = 0.1509078305790594
$ws.Cells.Item(4,10).Value = 0.1509078305790594
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 83.95844533333333
$ws.Cells.Item(4,14).Value = 251.875336
$ws.Cells.Item(4,15).Value = 0.997975448786732
$ws.Cells.Item(4,16).Value = 0.997975448786732
$ws.Cells.Item(4,17).Value = 2464.468320014678
$ws.Cells.Item(4,18).Value = 22180.2148801321
$ws.Cells.Item(4,19).Value = 0.1506023099475689
$ws.Cells.Item(4,20).Value = 0.1506023099475689

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Icam1"
$ws.Cells.Item(5,3).Value = "Itgb2"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 29.35342966666667
$ws.Cells.Item(5,8).Value = 88.060289
$ws.Cells.Item(5,9).Value = 0.1509078305790594
$ws.Cells.Item(5,10).Value = 0.1509078305790594
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.170323
$ws.Cells.Item(5,14).Value = 0.510969
$ws.Cells.Item(5,15).Value = 0.002024551213268089
$ws.Cells.Item(5,16).Value = 0.00202455121326809
$ws.Cells.Item(5,17).Value = 4.999564201115667
$ws.Cells.Item(5,18).Value = 44.996077810041
$ws.Cells.Item(5,19).Value = 0.00030552063149049
$ws.Cells.Item(5,20).Value = 0.00030552063149049

# Row 6
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Icam1"
$ws.Cells.Item(6,3).Value = "Itgb2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.248038666666667
$ws.Cells.Item(6,8).Value = 3.744116
$ws.Cells.Item(6,9).Value = 0.006416245386116614
$ws.Cells.Item(6,10).Value = 0.006416245386116614
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 83.95844533333333
$ws.Cells.Item(6,14).Value = 251.875336
$ws.Cells.Item(6,15).Value = 0.997975448786732
$ws.Cells.Item(6,16).Value = 0.997975448786732
$ws.Cells.Item(6,17).Value = 104.7833861692196
$ws.Cells.Item(6,18).Value = 943.050475522976
$ws.Cells.Item(6,19).Value = 0.006403255368735526
$ws.Cells.Item(6,20).Value = 0.006403255368735527

# Row 7
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Icam1"
$ws.Cells.Item(7,3).Value = "Itgb2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.248038666666667
$ws.Cells.Item(7,8).Value = 3.744116
$ws.Cells.Item(7,9).Value = 0.006416245386116614
$ws.Cells.Item(7,10).Value = 0.006416245386116614
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.170323
$ws.Cells.Item(7,14).Value = 0.510969
$ws.Cells.Item(7,15).Value = 0.002024551213268089
$ws.Cells.Item(7,16).Value = 0.00202455121326809
$ws.Cells.Item(7,17).Value = 0.2125696898226667
$ws.Cells.Item(7,18).Value = 1.913127208404
$ws.Cells.Item(7,19).Value = 0.00001299001738108817
$ws.Cells.Item(7,20).Value = 0.00001299001738108817

